$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: Transfer Volume changes from 250 to 500
$ws.Range("H2").Value = 500

# Add new row 3, duplicating most of row 2's content but with its own UID,
# destination well, and transfer volume (the original 250 that moved down).
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "level 2 384 source plate (DNA components)"
$ws.Range("C3").Value = "384LDV_AQ_B"
$ws.Range("D3").Value = "A3"
$ws.Range("E3").Value = "384-Well Level 2 MoClo output plate"
$ws.Range("F3").Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
$ws.Range("G3").Value = "A2"
$ws.Range("H3").Value = 250
$ws.Range("I3").Value = "pTU2-a-RFP"
